$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> @(PriceValue, VolumeValue). $null means "leave unchanged".
$changes = @{
    2  = @("43.159.66", "  +0.79%  ")
    3  = @("2.351.57", "  +4.99%  ")
    4  = @("1.01", "  -0.09%  ")
    5  = @("311.62", "  +5.37%  ")
    6  = @("108.88", "  -3.17%  ")
    7  = @("0.636", "  +1.40%  ")
    8  = @($null, "  -0.48%  ")
    9  = @($null, "  +1.01%  ")
    10 = @("43.60", "  -4.41%  ")
    11 = @("0.0939", "  +1.60%  ")
    12 = @("8.98", "  -0.11%  ")
    13 = @($null, "  +18.26%  ")
    14 = @($null, "  +0.74%  ")
    15 = @("16.31", "  +6.99%  ")
    16 = @("2.704.72", "  +4.77%  ")
    17 = @("2.412.22", "  +6.64%  ")
    18 = @("43.117.87", "  +0.96%  ")
    19 = @($null, "  +0.73%  ")
    20 = @("7.23", "  -3.59%  ")
    21 = @("75.72", "  +3.37%  ")
    22 = @("3.46", "  -2.40%  ")
    23 = @($null, "  +8.16%  ")
    24 = @("250.57", "  +8.34%  ")
    25 = @("8.99", "  -4.57%  ")
    26 = @("11.94", "  -0.42%  ")
    27 = @("1.00", "  +0.15%  ")
    28 = @("2.25", "  +0.72%  ")
    29 = @("38.74", "  -2.93%  ")
    30 = @("22.54", "  +6.05%  ")
    31 = @("174.27", "  +0.13%  ")
    33 = @($null, "  +1.63%  ")
    34 = @("5.81", "  +1.27%  ")
    35 = @("4.98", "  -1.15%  ")
    36 = @($null, "  +2.63%  ")
    37 = @($null, "  +2.02%  ")
    38 = @("4.14", "  -3.60%  ")
    39 = @($null, "  -1.13%  ")
    40 = @("2.76", "  +8.53%  ")
    41 = @($null, "  +12.79%  ")
    42 = @("72.24", "  +0.66%  ")
    43 = @("0.234", "  -2.24%  ")
    44 = @($null, "  -0.15%  ")
    45 = @("12.60", "  -4.92%  ")
    46 = @("5.71", "  +2.29%  ")
    47 = @("9.23", "  +6.58%  ")
    48 = @("110.51", "  +4.07%  ")
    49 = @($null, "  -1.51%  ")
    50 = @("0.100", "  +1.34%  ")
    51 = @("71.15", "  +4.68%  ")
}

foreach ($row in $changes.Keys) {
    $vals = $changes[$row]
    $priceVal = $vals[0]
    $volVal = $vals[1]

    if ($null -ne $priceVal) {
        $ws.Range("D$row").Value = $priceVal
    }
    if ($null -ne $volVal) {
        $ws.Range("E$row").Value = $volVal
    }
}
